# Refresh cryptos list (prices / 1h volume %) per GitHub Actions data pull.
# Note: D-column price cells that parse as plain decimals (e.g. "586.97")
# are written with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr cells, e.g. "64.244.90" which already
# fails numeric parsing because of the thousands-dot formatting used here).
# The cell Style is then reset to "Normal" so no stray quote-prefix / text
# number-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.244.90'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '3.491.02'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'586.97"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = "'134.07"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.32%  '
$ws.Range('D7').Value = '3.490.01'
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('D13').Value = '4.085.66'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '3.491.86'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = '64.298.69'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = "'25.24"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.54%  '
$ws.Range('E19').Value = '  -1.56%  '
$ws.Range('D20').Value = "'5.74"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').Value = "'13.60"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.15%  '
$ws.Range('D22').Value = "'387.79"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('D23').Value = '3.631.20'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('E24').Value = '  -2.51%  '
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = "'7.38"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('D33').Value = "'8.23"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').Value = '3.512.15'
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('D36').Value = "'0.149"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.72%  '
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('D38').Value = "'5.24"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('D41').Value = "'162.47"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.70%  '
$ws.Range('D42').Value = "'0.0782"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.34%  '
$ws.Range('E43').Value = '  -1.23%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = "'25.39"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.33%  '
$ws.Range('D46').Value = "'41.74"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'1.66"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = "'1.18"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.51%  '
$ws.Range('D50').Value = '2.471.57'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('D51').Value = "'6.75"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.19%  '
